$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raspberry Pi")

# Append the newly supported kernel version line to the Rpi4 / Rpi5 kernel
# version cells (row 4 of the "Raspberry Pi" sheet).
$ws.Range("C4").Value = "2023-12-05 6.1.0-rpi7-rpi-v8" + [char]10 + "2024-03-15 6.6.20+rpt-rpi-v8" + [char]10 + "2024-07-04 6.6.31+rpt-rpi-v8" + [char]10 + "2025-13-05 6.12.25+rpt-rpi-v8"
$ws.Range("G4").Value = "2024-07-04 6.6.31+rpt-rpi-2712" + [char]10 + "2025-13-05 6.12.25+rpt-rpi-2712"

# The extra line increases the amount of text wrapped in the row, so Excel
# grows row 4 to fit it.
$ws.Rows.Item(4).RowHeight = 75

# Update the active selection to match where the edit was made.
$ws.Range("G4").Select() | Out-Null
